# PanaudojimoAtvejuSpecifikacijos.docx — "Created sequence diagrams for theather."
#
# Three edits to word/document.xml, all inside the "PA 'Seansų užimtumo
# statistika'" use case table:
#
#   1. The run holding 'PA "Seansų užimtumo statistika"' gets split in two
#      (" “Seansų uži" | "mtumo statistika“ ") with the document's _GoBack
#      bookmark sitting at the split point (this is where Word's editor
#      cursor last stood).
#   2. The pre-condition cell text gains the missing word "atsidaręs":
#      "...administratorius ir įvykio langą." ->
#      "...administratorius ir atsidaręs įvykio langą."
#   3. The _GoBack bookmark that used to sit mid-word in "u|žpildyta" (right
#      after "p") is removed, since _GoBack moved to the new location in #1.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: delete the (single) bookmark whose zero-length anchor sits at
# document position $pos, leaving the surrounding text untouched. A
# zero-width bookmark is only dropped by the engine when the deleted
# range strictly straddles its anchor, so we snapshot the one character
# on each side, delete that 2-char span (which also removes the
# bookmark that lived in between) and type the same two characters back.
# ---------------------------------------------------------------------
function Remove-BookmarkAt($pos) {
    $span = $d.Range($pos - 1, $pos + 1)
    $savedText = $span.Text
    $span.Delete()
    $restore = $d.Range($pos - 1, $pos - 1)
    $restore.InsertBefore($savedText)
}

# ---------------------------------------------------------------------
# 1 & 3. Relocate the _GoBack bookmark from its old spot (inside
#         "u žpildyta", right after the lone "p" run) to the new spot
#         (inside "Seansų užimtumo statistika", right after "uži").
# ---------------------------------------------------------------------

# Old location: the run sequence ...">p</w:t></w:r><w:bookmarkStart .../>
# <w:bookmarkEnd .../><w:r ...>"ildyta. ...". Find the unique "p" + "ildyta"
# boundary via the surrounding unique text.
$oldFound = $d.Content.Duplicate
$oldFound.Find.Execute("salėse užp")
$oldPos = $oldFound.End   # right after the lone "p" run, i.e. where the bookmark sits

# New location: inside "Seansų užimtumo statistika", right after "uži".
$newFound = $d.Content.Duplicate
$newFound.Find.Execute("Seansų uži")
$newPos = $newFound.End

# Add the bookmark at the new position first (names can repeat here; we
# remove the stray duplicate immediately after), then drop the old one.
$newRange = $d.Range($newPos, $newPos)
$d.Bookmarks.Add("_GoBack", $newRange)
Remove-BookmarkAt $oldPos

# ---------------------------------------------------------------------
# 2. Insert the missing word "atsidaręs" into the pre-condition sentence
#    for the "Seansų užimtumo statistika" use case.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Vartotojas turi būti prisijungęs kaip kino teatro administratorius ir įvykio langą.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Vartotojas turi būti prisijungęs kaip kino teatro administratorius ir atsidaręs įvykio langą.",
    2)
